$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Oct 06 11:23:15 EDT 2023"
$ws.Range("B3").Value = "Fri Oct 06 11:23:28 EDT 2023"
$ws.Range("B4").Value = "Fri Oct 06 11:23:41 EDT 2023"
$ws.Range("B5").Value = "Fri Oct 06 11:23:54 EDT 2023"
$ws.Range("B6").Value = "Fri Oct 06 11:24:07 EDT 2023"
$ws.Range("B7").Value = "Fri Oct 06 11:24:20 EDT 2023"
$ws.Range("B8").Value = "Fri Oct 06 11:24:33 EDT 2023"
$ws.Range("B9").Value = "Fri Oct 06 11:24:47 EDT 2023"
$ws.Range("B10").Value = "Fri Oct 06 11:25:00 EDT 2023"
$ws.Range("B11").Value = "Fri Oct 06 11:25:13 EDT 2023"
$ws.Range("B12").Value = "Fri Oct 06 11:25:26 EDT 2023"
$ws.Range("B13").Value = "Fri Oct 06 11:25:39 EDT 2023"
$ws.Range("B14").Value = "Fri Oct 06 11:25:52 EDT 2023"
$ws.Range("B15").Value = "Fri Oct 06 11:26:05 EDT 2023"
$ws.Range("B16").Value = "Fri Oct 06 11:26:18 EDT 2023"
$ws.Range("B17").Value = "Fri Oct 06 11:26:31 EDT 2023"
$ws.Range("B18").Value = "Fri Oct 06 11:26:44 EDT 2023"
$ws.Range("B19").Value = "Fri Oct 06 11:26:57 EDT 2023"
$ws.Range("B20").Value = "Fri Oct 06 11:27:10 EDT 2023"
$ws.Range("B21").Value = "Fri Oct 06 11:27:23 EDT 2023"
$ws.Range("B22").Value = "Fri Oct 06 11:27:36 EDT 2023"
$ws.Range("B23").Value = "Fri Oct 06 11:27:49 EDT 2023"
$ws.Range("B24").Value = "Fri Oct 06 11:28:02 EDT 2023"
$ws.Range("B25").Value = "Fri Oct 06 11:28:15 EDT 2023"
$ws.Range("B26").Value = "Fri Oct 06 11:28:28 EDT 2023"
$ws.Range("B27").Value = "Fri Oct 06 11:28:41 EDT 2023"
$ws.Range("B28").Value = "Fri Oct 06 11:28:54 EDT 2023"
$ws.Range("B29").Value = "Fri Oct 06 11:29:07 EDT 2023"
$ws.Range("B30").Value = "Fri Oct 06 11:29:20 EDT 2023"
$ws.Range("B31").Value = "Fri Oct 06 11:29:33 EDT 2023"
$ws.Range("B32").Value = "Fri Oct 06 11:29:46 EDT 2023"
$ws.Range("B33").Value = "Fri Oct 06 11:30:00 EDT 2023"
